# Daily attendance processing - 2025-12-14 01:36:43
# Reorder the "Recorded By" (column G) contributor list on each attendance
# row so that any "System" / "system" entries sort after the real (human /
# integration) recorder addresses, instead of leading the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 157
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ", "
    $nonSystem = @()
    $systemParts = @()
    foreach ($p in $parts) {
        if ($p.Trim().ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $nonSystem += $p
        }
    }

    $newParts = $nonSystem + $systemParts
    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
